# Update the Cal Poly_A team-specific transition matrix sheet
# with newly collected time-data (row-normalized probabilities).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1868512110726644
$ws.Range("C2").Value = 0.5570934256055363
$ws.Range("J2").Value = 0.006920415224913495
$ws.Range("P2").Value = 0.1280276816608996
$ws.Range("S2").Value = 0.1211072664359862
$ws.Range("B3").Value = 0.006097560975609756
$ws.Range("C3").Value = 0.01219512195121951
$ws.Range("J3").Value = 0.02439024390243903
$ws.Range("P3").Value = 0.7317073170731707
$ws.Range("S3").Value = 0.225609756097561
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.6046511627906976
$ws.Range("S4").Value = 0.3720930232558139
$ws.Range("B6").Value = 0.09124087591240876
$ws.Range("D6").Value = 0.0072992700729927
$ws.Range("F6").Value = 0.06204379562043796
$ws.Range("J6").Value = 0.218978102189781
$ws.Range("O6").Value = 0.0218978102189781
$ws.Range("Q6").Value = 0.1240875912408759
$ws.Range("R6").Value = 0.08029197080291971
$ws.Range("S6").Value = 0.3941605839416059
$ws.Range("B7").Value = 0.1061224489795918
$ws.Range("D7").Value = 0.02857142857142857
$ws.Range("F7").Value = 0.07755102040816327
$ws.Range("J7").Value = 0.1265306122448979
$ws.Range("O7").Value = 0.01224489795918367
$ws.Range("Q7").Value = 0.1551020408163265
$ws.Range("R7").Value = 0.1020408163265306
$ws.Range("S7").Value = 0.3918367346938775
$ws.Range("B8").Value = 0.07569721115537849
$ws.Range("D8").Value = 0.02390438247011952
$ws.Range("E8").Value = 0.00199203187250996
$ws.Range("F8").Value = 0.06772908366533864
$ws.Range("J8").Value = 0.099601593625498
$ws.Range("O8").Value = 0.02589641434262948
$ws.Range("Q8").Value = 0.1912350597609562
$ws.Range("R8").Value = 0.05577689243027888
$ws.Range("S8").Value = 0.4581673306772908
$ws.Range("B9").Value = 0.07975460122699386
$ws.Range("D9").Value = 0.006134969325153374
$ws.Range("F9").Value = 0.09202453987730061
$ws.Range("J9").Value = 0.09815950920245399
$ws.Range("O9").Value = 0.01226993865030675
$ws.Range("Q9").Value = 0.1595092024539877
$ws.Range("R9").Value = 0.09815950920245399
$ws.Range("S9").Value = 0.4539877300613497
$ws.Range("B10").Value = 0.1057304277643261
$ws.Range("D10").Value = 0.01937046004842615
$ws.Range("F10").Value = 0.08716707021791767
$ws.Range("J10").Value = 0.1138014527845036
$ws.Range("O10").Value = 0.02502017756255044
$ws.Range("Q10").Value = 0.1799838579499596
$ws.Range("R10").Value = 0.06456820016142049
$ws.Range("S10").Value = 0.4043583535108959
$ws.Range("G11").Value = 0.1272727272727273
$ws.Range("J11").Value = 0.1012987012987013
$ws.Range("K11").Value = 0.1688311688311688
$ws.Range("L11").Value = 0.587012987012987
$ws.Range("S11").Value = 0.01558441558441558
$ws.Range("G12").Value = 0.7264957264957265
$ws.Range("J12").Value = 0.2051282051282051
$ws.Range("K12").Value = 0.008547008547008548
$ws.Range("L12").Value = 0.0170940170940171
$ws.Range("S12").Value = 0.04273504273504274
$ws.Range("G13").Value = 0.7254901960784313
$ws.Range("J13").Value = 0.2549019607843137
$ws.Range("S13").Value = 0.0196078431372549
$ws.Range("F15").Value = 0.01834862385321101
$ws.Range("H15").Value = 0.1880733944954129
$ws.Range("I15").Value = 0.03669724770642202
$ws.Range("J15").Value = 0.3302752293577982
$ws.Range("K15").Value = 0.05045871559633028
$ws.Range("M15").Value = 0.01834862385321101
$ws.Range("N15").Value = 0.009174311926605505
$ws.Range("O15").Value = 0.02752293577981652
$ws.Range("S15").Value = 0.3211009174311927
$ws.Range("F16").Value = 0.01675977653631285
$ws.Range("H16").Value = 0.2067039106145251
$ws.Range("I16").Value = 0.0335195530726257
$ws.Range("J16").Value = 0.3743016759776536
$ws.Range("K16").Value = 0.1284916201117318
$ws.Range("M16").Value = 0.01675977653631285
$ws.Range("O16").Value = 0.02793296089385475
$ws.Range("S16").Value = 0.1955307262569832
$ws.Range("F17").Value = 0.01927710843373494
$ws.Range("H17").Value = 0.1951807228915663
$ws.Range("I17").Value = 0.0891566265060241
$ws.Range("J17").Value = 0.380722891566265
$ws.Range("K17").Value = 0.1349397590361446
$ws.Range("M17").Value = 0.02409638554216868
$ws.Range("N17").Value = 0.002409638554216868
$ws.Range("O17").Value = 0.05301204819277108
$ws.Range("S17").Value = 0.1012048192771084
$ws.Range("F18").Value = 0.01775147928994083
$ws.Range("H18").Value = 0.1952662721893491
$ws.Range("I18").Value = 0.08284023668639054
$ws.Range("J18").Value = 0.378698224852071
$ws.Range("K18").Value = 0.1420118343195266
$ws.Range("M18").Value = 0.005917159763313609
$ws.Range("N18").Value = 0.005917159763313609
$ws.Range("O18").Value = 0.04142011834319527
$ws.Range("S18").Value = 0.1301775147928994
$ws.Range("F19").Value = 0.02126200274348422
$ws.Range("H19").Value = 0.2139917695473251
$ws.Range("I19").Value = 0.06858710562414266
$ws.Range("J19").Value = 0.3381344307270233
$ws.Range("K19").Value = 0.135116598079561
$ws.Range("M19").Value = 0.02400548696844993
$ws.Range("O19").Value = 0.06790123456790123
$ws.Range("S19").Value = 0.1310013717421125

Write-Host "Applied updated transition probabilities."
